$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Quarterly labels that replace the serial-date values in column A (rows 2-22).
$labels = @(
    "2004Q4", "2005Q4", "2006Q4", "2007Q4", "2008Q4", "2009Q4", "2010Q4",
    "2011Q4", "2012Q4", "2013Q4", "2014Q4", "2015Q4", "2016Q4", "2017Q4",
    "2018Q4", "2019Q4", "2020Q4", "2021Q4", "2022Q4", "2023Q4", "2024Q4"
)

for ($i = 0; $i -lt $labels.Length; $i++) {
    $row = $i + 2
    $cell = $ws.Cells.Item($row, 1)
    $cell.Value = $labels[$i]
}

# Re-apply the same (non date-formatted) style used by the header row so the
# now-obsolete custom date number format is dropped from the workbook.
$ws.Range("A1").Copy()
$ws.Range("A2:A22").PasteSpecial(-4122) # xlPasteFormats
$excel.CutCopyMode = $false
